$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.868.21'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.318.88'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'302.39"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = "'96.23"
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").Value = "'0.507"
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = "'34.57"
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").Value = "'18.93"
$ws.Range("E11").Value = '  +5.51%  '
$ws.Range("D12").Value = "'0.0785"
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = "'6.77"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '2.684.40'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = '2.313.17'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").Value = '42.809.26'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = "'12.19"
$ws.Range("E19").Value = '  -6.20%  '
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").Value = '0.0₃0893'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = "'67.95"
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = '  +6.63%  '
$ws.Range("D24").Value = "'236.30"
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = "'24.40"
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("D31").Value = "'32.30"
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").Value = "'17.88"
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").Value = "'0.0701"
$ws.Range("E36").Value = '  +2.62%  '
$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = '  +3.17%  '
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = '  +3.37%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").Value = "'20.73"
$ws.Range("E42").Value = '  +13.11%  '
$ws.Range("D43").Value = '1.935.72'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").Value = "'10.20"
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").Value = '  +3.55%  '
$ws.Range("D47").Value = "'2.76"
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = '2.551.58'
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").Value = "'53.42"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").Value = "'2.83"
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").Value = "'72.12"
$ws.Range("E51").Value = '  +2.47%  '
